# Add a new slide ("How the algorithm will work:") as slide 5,
# using the same "Title and Content" layout as the other content slides.

$p = $ppt.ActivePresentation

# ppLayoutText = 2 -> Title + Content placeholder layout (matches slideLayout2.xml
# already used by slides 2-4 in this deck).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# ---------------------------------------------------------------------------
# Title placeholder
# ---------------------------------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "How the algorithm will work:"
# Matches <a:ea typeface="Calibri Light"/><a:cs typeface="Calibri Light"/>
$title.Font.NameFarEast = "Calibri Light"
$title.Font.NameComplexScript = "Calibri Light"

# ---------------------------------------------------------------------------
# Body / content placeholder
# ---------------------------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame

# Body paragraph text, in order. Paragraph 1 sits at the top outline level;
# the rest are indented one level (lvl=1 / IndentLevel=2).
$paraTexts = @(
    "The algorithm will follow a greedy approach:",
    "Start with an empty schedule.",
    "Loop the schedule available slot.",
    "For each slot pick all the possible availabilities that might fit the slot and not been selected yet.",
    "Sort those availabilities according to priority factors, like prioritize the required courses over elective courses.",
    "Loop those availabilities from the beginning and pick what ever availability that does not made a conflict in the professor schedule or year schedule.",
    "Assign this availability to the time slot and mark it as selected",
    ""
)

# Build up the text one paragraph at a time. Re-assigning the whole
# TextRange.Text (instead of InsertAfter) on every step makes each newly
# split-off paragraph clone the east-asian/complex-script typeface that was
# already applied to paragraph 1 -- this is what lets every paragraph end up
# with <a:ea typeface="Calibri"/><a:cs typeface="Calibri"/>.
$full = $paraTexts[0]
$body.TextRange.Text = $full
$body.TextRange.Font.NameFarEast = "Calibri"
$body.TextRange.Font.NameComplexScript = "Calibri"

for ($i = 1; $i -lt $paraTexts.Count; $i++) {
    $full = $full + "`r" + $paraTexts[$i]
    $body.TextRange.Text = $full
}

# Outline levels: paragraph 1 stays at the top level, paragraphs 2-8 are
# indented one level (OOXML lvl="1" == COM IndentLevel 2).
for ($i = 2; $i -le $paraTexts.Count; $i++) {
    $body.TextRange.Paragraphs($i, 1).IndentLevel = 2
}

# Bold "greedy " inside paragraph 1.
$boldStart = $paraTexts[0].IndexOf("greedy ") + 1
$body.TextRange.Characters($boldStart, 7).Font.Bold = $true

# Body text box formatting: vert="horz", 0.1"/0.05" insets, anchored top,
# shrink text on overflow (<a:normAutofit/>).
$body.Orientation = 1
$body.MarginLeft = 7.2
$body.MarginTop = 3.6
$body.MarginRight = 7.2
$body.MarginBottom = 3.6
$body.VerticalAnchor = 1
$body.AutoSize = 2
